# Refresh the cryptos price/volume snapshot (and fix the two swapped/replaced rows)
# to match the source feed pulled on Tue Dec 12 07:40:02 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe forces Excel to keep numeric-looking "Price" strings as
# plain text (matching the original cells, which are all stored as text) instead
# of silently converting them to numbers.
$apos = "'"

$ws.Range("D2").Value = $apos + '41.866.50'
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").Value = $apos + '2.226.49'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = $apos + '250.98'
$ws.Range("E5").Value = '  +7.87%  '
$ws.Range("D6").Value = $apos + '0.631'
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("D7").Value = $apos + '71.40'
$ws.Range("E7").Value = '  +2.79%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = $apos + '0.590'
$ws.Range("E9").Value = '  +5.29%  '
$ws.Range("D10").Value = $apos + '41.35'
$ws.Range("E10").Value = '  +15.10%  '
$ws.Range("D11").Value = $apos + '0.0969'
$ws.Range("E11").Value = '  -2.01%  '
$ws.Range("D12").Value = $apos + '58.23'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").Value = $apos + '7.14'
$ws.Range("E13").Value = '  +5.87%  '
$ws.Range("D14").Value = $apos + '0.106'
$ws.Range("E14").Value = '  +1.14%  '
$ws.Range("D15").Value = $apos + '2.558.82'
$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("D16").Value = $apos + '14.95'
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("E17").Value = '  +1.37%  '
$ws.Range("D18").Value = $apos + '2.227.49'
$ws.Range("E18").Value = '  -0.44%  '
$ws.Range("D19").Value = $apos + '41.762.98'
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("D22").Value = $apos + '73.01'
$ws.Range("E22").Value = '  -0.40%  '
$ws.Range("D23").Value = $apos + '235.30'
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("E24").Value = '  +7.80%  '
$ws.Range("D25").Value = $apos + '4.19'
$ws.Range("E25").Value = '  +14.45%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").Value = $apos + '2.54'
$ws.Range("E27").Value = '  +7.67%  '
$ws.Range("E28").Value = '  +6.53%  '
$ws.Range("E29").Value = '  +0.52%  '
$ws.Range("D30").Value = $apos + '170.93'
$ws.Range("E30").Value = '  +1.13%  '
$ws.Range("D31").Value = $apos + '20.76'
$ws.Range("E31").Value = '  +0.65%  '
$ws.Range("D32").Value = $apos + '0.122'
$ws.Range("E32").Value = '  +2.74%  '
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").Value = $apos + '0.125'
$ws.Range("E33").Value = '  -1.35%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = $apos + '5.60'
$ws.Range("E34").Value = '  +6.09%  '
$ws.Range("D35").Value = $apos + '0.0725'
$ws.Range("E35").Value = '  +1.70%  '
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("D37").Value = $apos + '26.11'
$ws.Range("E37").Value = '  +19.49%  '
$ws.Range("D38").Value = $apos + '3.92'
$ws.Range("E38").Value = '  +8.39%  '
$ws.Range("D39").Value = $apos + '0.0302'
$ws.Range("E39").Value = '  +13.73%  '
$ws.Range("E40").Value = '  +2.41%  '
$ws.Range("E41").Value = '  -0.81%  '
$ws.Range("D42").Value = $apos + '68.00'
$ws.Range("E42").Value = '  +1.42%  '
$ws.Range("E43").Value = '  +9.60%  '
$ws.Range("D44").Value = $apos + '11.74'
$ws.Range("E44").Value = '  +16.74%  '
$ws.Range("D45").Value = $apos + '8.83'
$ws.Range("D46").Value = $apos + '4.81'
$ws.Range("E46").Value = '  -2.84%  '
$ws.Range("E47").Value = '  +1.86%  '
$ws.Range("D48").Value = $apos + '4.72'
$ws.Range("E48").Value = '  +8.16%  '
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("E50").Value = '  +7.94%  '
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").Value = $apos + '1.19'
$ws.Range("E51").Value = '  +1.27%  '
